# Auto-generated edit script: apply cryptocurrency price/volume updates
# per commit 'Updated symbol list on Wed Feb 15 14:42:12 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '303.01'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.18%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '43.30'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '7.65%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.052'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.64%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07668'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.53%'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.407'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.22%'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.605'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.70%'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.006'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '8.91%'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1244'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '6.88%'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1857'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '4.39%'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09099'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.88%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04170'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.47%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1046'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.44%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001271'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.35%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005755'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-2.38%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1,893.71%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.330'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.47%'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.356'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.93%'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3353'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.00%'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.444'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '6.32%'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1397'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.11%'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.3198'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '13.90%'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04165'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '5.14%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001283'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.97%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004468'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '17.54%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '9.41%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02452'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '5.10%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05284'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '4.06%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.005967'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.33%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007688'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.06%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1349'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4.29%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.75%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007462'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.92%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3027'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '3.78%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006702'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '8.12%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.42%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.04094'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-11.48%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.12%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002096'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.42%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001996'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.42%'
